$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "Cython " -> split into a spell-checked "Cython" run + a separate
#    trailing-space run (proofErr wraps the word that the spell-checker
#    flagged).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Cython -> Bridging between python and C, probably would help a bit", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find Cython paragraph" }
$para = $rng.Paragraphs(1)
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Cython</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>-&gt; Bridging between python and C, probably would help a bit</w:t>
  </w:r>
</w:p>
"@
$para.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) "...Reddit, Youtube, and other sources..." -> split "Youtube" into its
#    own spell-checked run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("because just feeding them past costs isn", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find ML Models paragraph" }
$para = $rng.Paragraphs(1)
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">Reinforcement? ML Model -&gt; </w:t>
  </w:r>
  <w:r>
    <w:t>ML Models suck at predicting future costs</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">because just feeding them past costs isn&#8217;t the full picture. Reddit, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Youtube</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>, and other sources would provide a much truer picture.</w:t>
  </w:r>
</w:p>
"@
$para.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3) NEW GRAPH paragraph: split "post_id" into its own spell-checked run,
#    then collapse the run of 8 empty placeholder paragraphs down to just
#    2 blanks, with a brand new bullet "Run algorithm on graphics card?"
#    in between (the GPU idea from the commit message).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("NEW GRAPH:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find NEW GRAPH paragraph" }
$newGraphPara = $rng.Paragraphs(1)

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Monetization:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find Monetization paragraph" }
$monetizationPara = $rng2.Paragraphs(1)

$blockRange = $d.Range($newGraphPara.Range.Start, $monetizationPara.Range.Start)
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">NEW GRAPH: </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Relation in DB that tracks how any views/votes/comments a stock gets per day, regardless of </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>post_id</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>, CREATE TABLE (symbol, date, votes, views, comments), PRIMARY KEY(symbol, date)</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Run algorithm on graphics card?</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
</w:p>
"@
$blockRange.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4) "Monetization:" title gains a <w:lastRenderedPageBreak/> before its
#    text run (the new content above pushed it onto a fresh page).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Monetization:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find Monetization paragraph (pass 2)" }
$para = $rng.Paragraphs(1)
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="Title"/>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>Monetization:</w:t>
  </w:r>
</w:p>
"@
$para.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 5) "Sell interactive data_visualization tool" -> split
#    "data_visualization" into its own spell-checked run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Sell interactive data_visualization tool", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find Sell interactive paragraph" }
$para = $rng.Paragraphs(1)
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">Sell interactive </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>data_visualization</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> tool</w:t>
  </w:r>
</w:p>
"@
$para.Range.InsertXML($xml)

Write-Output "edits applied"
